# Atualização de bases das ligas, do dia: 03-05-2024 às 22:15
#
# The source feed re-sorted a handful of fixtures that share the same
# kickoff date, which swapped the data rows for the following pairs
# (the leading "id"/index column A keeps the row's position, only the
# match data in columns B:AB move):
#   rows 148 <-> 150
#   rows 153 <-> 154
#   rows 211 <-> 212
#   rows 214 <-> 215

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData {
    param($RowA, $RowB)

    $rangeA = $ws.Range("B$RowA" + ":AB$RowA")
    $rangeB = $ws.Range("B$RowB" + ":AB$RowB")

    $valuesA = $rangeA.Value2
    $valuesB = $rangeB.Value2

    $rangeA.Value2 = $valuesB
    $rangeB.Value2 = $valuesA
}

Swap-RowData 148 150
Swap-RowData 153 154
Swap-RowData 211 212
Swap-RowData 214 215
